# Applies the commit "Refined metadata to be additional tab":
#   1. Updates the "data" sheet's time_taken (column F) timestamps for rows 2..93.
#   2. Adds a new "metadata" worksheet after "data", summarising the query
#      that produced the "data" sheet's contents.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1) Update the time_taken (column F) timestamps on the "data" sheet.
# ---------------------------------------------------------------------------
$newTimes = @{
    2 = "2021-10-05 14:34:02.702105"
    3 = "2021-10-05 14:34:02.702113"
    4 = "2021-10-05 14:34:02.702116"
    5 = "2021-10-05 14:34:02.702118"
    6 = "2021-10-05 14:34:02.702121"
    7 = "2021-10-05 14:34:02.702124"
    8 = "2021-10-05 14:34:02.702126"
    9 = "2021-10-05 14:34:02.702129"
    10 = "2021-10-05 14:34:02.702132"
    11 = "2021-10-05 14:34:02.702134"
    12 = "2021-10-05 14:34:02.702137"
    13 = "2021-10-05 14:34:02.702139"
    14 = "2021-10-05 14:34:02.702142"
    15 = "2021-10-05 14:34:02.702144"
    16 = "2021-10-05 14:34:02.702147"
    17 = "2021-10-05 14:34:02.702149"
    18 = "2021-10-05 14:34:02.702152"
    19 = "2021-10-05 14:34:02.702155"
    20 = "2021-10-05 14:34:02.702157"
    21 = "2021-10-05 14:34:02.702160"
    22 = "2021-10-05 14:34:02.702162"
    23 = "2021-10-05 14:34:02.702165"
    24 = "2021-10-05 14:34:02.702167"
    25 = "2021-10-05 14:34:02.702169"
    26 = "2021-10-05 14:34:02.702172"
    27 = "2021-10-05 14:34:02.702175"
    28 = "2021-10-05 14:34:02.702177"
    29 = "2021-10-05 14:34:02.702180"
    30 = "2021-10-05 14:34:02.702182"
    31 = "2021-10-05 14:34:02.702185"
    32 = "2021-10-05 14:34:02.702187"
    33 = "2021-10-05 14:34:02.702190"
    34 = "2021-10-05 14:34:02.702193"
    35 = "2021-10-05 14:34:02.702195"
    36 = "2021-10-05 14:34:02.702198"
    37 = "2021-10-05 14:34:02.702200"
    38 = "2021-10-05 14:34:02.702203"
    39 = "2021-10-05 14:34:02.702205"
    40 = "2021-10-05 14:34:02.702208"
    41 = "2021-10-05 14:34:02.702210"
    42 = "2021-10-05 14:34:02.702213"
    43 = "2021-10-05 14:34:02.702216"
    44 = "2021-10-05 14:34:02.702218"
    45 = "2021-10-05 14:34:02.702221"
    46 = "2021-10-05 14:34:02.702223"
    47 = "2021-10-05 14:34:02.702226"
    48 = "2021-10-05 14:34:02.702228"
    49 = "2021-10-05 14:34:02.702230"
    50 = "2021-10-05 14:34:02.702233"
    51 = "2021-10-05 14:34:02.702235"
    52 = "2021-10-05 14:34:02.702238"
    53 = "2021-10-05 14:34:02.702240"
    54 = "2021-10-05 14:34:02.702243"
    55 = "2021-10-05 14:34:02.702246"
    56 = "2021-10-05 14:34:02.702248"
    57 = "2021-10-05 14:34:02.702251"
    58 = "2021-10-05 14:34:02.702253"
    59 = "2021-10-05 14:34:02.702256"
    60 = "2021-10-05 14:34:02.702258"
    61 = "2021-10-05 14:34:02.702261"
    62 = "2021-10-05 14:34:02.702263"
    63 = "2021-10-05 14:34:02.702266"
    64 = "2021-10-05 14:34:02.702268"
    65 = "2021-10-05 14:34:02.702271"
    66 = "2021-10-05 14:34:02.702274"
    67 = "2021-10-05 14:34:02.702277"
    68 = "2021-10-05 14:34:02.702280"
    69 = "2021-10-05 14:34:02.702282"
    70 = "2021-10-05 14:34:02.702285"
    71 = "2021-10-05 14:34:02.702288"
    72 = "2021-10-05 14:34:02.702290"
    73 = "2021-10-05 14:34:02.702293"
    74 = "2021-10-05 14:34:02.702295"
    75 = "2021-10-05 14:34:02.702298"
    76 = "2021-10-05 14:34:02.702300"
    77 = "2021-10-05 14:34:02.702303"
    78 = "2021-10-05 14:34:02.702308"
    79 = "2021-10-05 14:34:02.702311"
    80 = "2021-10-05 14:34:02.702313"
    81 = "2021-10-05 14:34:02.702316"
    82 = "2021-10-05 14:34:02.702318"
    83 = "2021-10-05 14:34:02.702321"
    84 = "2021-10-05 14:34:02.702323"
    85 = "2021-10-05 14:34:02.702326"
    86 = "2021-10-05 14:34:02.702328"
    87 = "2021-10-05 14:34:02.702330"
    88 = "2021-10-05 14:34:02.702333"
    89 = "2021-10-05 14:34:02.702336"
    90 = "2021-10-05 14:34:02.702338"
    91 = "2021-10-05 14:34:02.702340"
    92 = "2021-10-05 14:34:02.702343"
    93 = "2021-10-05 14:34:02.702345"
}
foreach ($row in $newTimes.Keys) {
    $dataSheet.Cells.Item([int]$row, 6).Value = $newTimes[$row]
}

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" worksheet right after "data".
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row: values first, then clone "data"'s header style (bold/border/
# centered) onto B1:G1 so the new sheet reuses the existing style index
# instead of registering a near-duplicate one.
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row (row 2). A2 reuses the same index-column style as data!A2.
$metaSheet.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$metaSheet.Range("B2").Value = "Hydrocephalus_Ventriculomegaly"
$metaSheet.Range("C2").Value = 115

# data_version ("0.97") must stay textual, not become the number 0.97.
$versionCell = $metaSheet.Range("D2")
$versionCell.NumberFormat = "@"
$versionCell.Value = "0.97"
$versionCell.Style = "Normal"

$metaSheet.Range("E2").Value = "2021-10-04T04:43:25.926555Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:34:02.698423"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/115/?format=json"

$excel.CutCopyMode = $false

# Restore focus to the "data" sheet (matches original activeTab=0).
$dataSheet.Activate()

Write-Host "metadata sheet added; data!F2:F93 timestamps updated"
